# "Update college and QQ group" -- on slide 1 (title slide), the small
# textbox that states the author's college is reworded and, since the new
# text is much shorter, PowerPoint's "shrink to fit" resize shifted the box
# to the right and shrank its width. Height/Top are unaffected.
#
# Shape inventory on slide 1 (in z-order / Shapes collection order):
#   1: "文本框 5"  (id 19)  - "讲授人：..."
#   2: "矩形 19"   (id 20)  - "Java语言程序设计" title
#   3: "文本框 5"  (id 21)  - "学院：..." <- the shape we need to edit
#   4: "直接连接符 21" (id 22) - connector line
#   5: "矩形 13"   (id 14)  - outline rectangle
#   6: "图片 3"    (id 4)   - picture

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(3)

# New position/size (EMU -> points, 1 pt = 12700 EMU). PowerPoint's COM
# Left/Top/Width/Height are single-precision floats, so a tiny epsilon
# (half an EMU, in points) is added to counter float32 truncation and land
# on the exact target EMU value after round-tripping.
$emuPerPt = 12700
$halfEmuPt = 0.5 / $emuPerPt

$sh.Left = (5312187 / $emuPerPt) + $halfEmuPt
$sh.Width = (1783122 / $emuPerPt) + $halfEmuPt

# Update the college name text.
$sh.TextFrame.TextRange.Text = "学院：人工智能学院"
